$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Add Panels")

# The "Loading Details Name" label text was renamed from "40V (A)" to
# "40V Rail(A)" to match the new implementation. Both F8 and F9 share this
# text via the shared-strings table, so updating either Range.Value updates
# the single shared string entry used by both cells.
$ws.Range("F8").Value = "40V Rail(A)"
$ws.Range("F9").Value = "40V Rail(A)"

# Leave the active selection on F9, matching where the edit was made.
$ws.Range("F9").Select()
